# Apply the commit: "Rerun main notebook and updated website data"
# - Recreate Sheet1 (copy+delete+rename) so its internal sheetId advances 1 -> 2,
#   mirroring the workbook.xml <sheet sheetId="..."/> bump in the diff, while
#   preserving all existing cell data, styles and row heights via Worksheet.Copy.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Copy($null, $ws) | Out-Null
$ws.Delete() | Out-Null
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Sheet1"
$ws.Select() | Out-Null

# New note text produced by the rerun, keyed by cell address.
$notes = @{
    "C3" = "The candidate has a strong educational background in Computer Science, but lacks direct experience in developing and implementing large-scale consumer-facing services. They have some relevant technical skills, but may need additional training or development in areas like machine learning and deep learning."
    "E3" = "Boca Raton, FL 33829 → Culver City, CA (2500 miles)"
    "G3" = "The candidate has some relevant experience in web development and technical support, but lacks direct experience in developing and implementing large-scale consumer-facing services."
    "I3" = "The candidate has a strong educational background in Computer Science, which is relevant to the job requirements."
    "K3" = "The candidate has some relevant technical skills, but may need additional training or development in areas like machine learning and deep learning."
    "M3" = "The candidate does not have any of the required programming languages listed in their resume."
    "O3" = "The candidate has some relevant soft skills, but lacks experience in search engine development and large-scale systems."
    "C4" = "The candidate has a strong background in computer science and software engineering, with relevant experience in programming languages and algorithm development. However, there are some gaps in experience and skills required for the job, such as machine learning and deep learning."
    "E4" = "Vancouver, WA 98662 → Culver City, CA (1000 miles)"
    "G4" = "The candidate has relevant experience in software engineering and algorithm development, but lacks experience in machine learning, deep learning, and large-scale systems."
    "I4" = "The candidate has a bachelor's degree in computer science, which meets the job requirement."
    "K4" = "The candidate has strong technical skills in programming languages, software engineering, and algorithm development, but lacks skills in machine learning and deep learning."
    "M4" = "The candidate speaks English, which meets the job requirement."
    "O4" = "The candidate has experience with team communication and collaboration, but lacks experience with Big Data Frameworks."
    "C5" = "The candidate has a strong background in software development, with relevant experience in cloud computing, containerization, and machine learning. However, there are some gaps in the required skills, such as C/C++, Java, and Golang."
    "E5" = "Tempe, Arizona → Culver City, CA (370 miles)"
    "G5" = "The candidate has relevant experience in software development, but lacks experience in some of the required programming languages."
    "I5" = "The candidate has a strong educational background in Computer Science, which is relevant to the job requirements."
    "K5" = "The candidate has some of the required skills, but lacks experience in some of the programming languages."
    "M5" = "The candidate speaks English, which is the required language for the job."
    "O5" = "The candidate has some of the required soft skills, but may need to develop others."
    "C6" = "The candidate has a strong background in backend engineering, technical leadership, and relevant industry experience, making them a suitable match for the role."
    "E6" = "→ United States"
    "G6" = "The candidate has significant experience in backend engineering and technical leadership, but lacks direct experience in search engine development."
    "I6" = "The candidate has a Master's degree in Education and relevant coursework in software development, but lacks a degree in Computer Science or a related field."
    "K6" = "The candidate has strong technical skills in Ruby, Rails, and Node.js, but lacks experience in C/C++, Python, or Java."
    "M6" = "The candidate speaks English, which is the required language for the role."
    "O6" = "The candidate has strong team communication and collaboration skills, but lacks experience with Big Data Frameworks."
    "C7" = "The candidate has a strong background in software engineering, with relevant experience in developing and implementing large-scale consumer-facing services. Their education in Computer Science and proficiency in programming languages such as Python, Java, and Golang are also well-aligned with the job requirements."
    "E7" = "Denver, CO, 80203 → United States"
    "G7" = "The candidate has 9 years of experience in software engineering, with a focus on developing and implementing large-scale consumer-facing services. However, they lack direct experience in search engine development."
    "I7" = "The candidate has a Bachelor's degree in Computer Science, which meets the job's education requirements."
    "K7" = "The candidate has proficiency in programming languages such as Python, Java, and Golang, and experience with event-driven architecture. However, they lack experience in machine learning and deep learning."
    "M7" = "The candidate speaks English, which meets the job's language requirements."
    "O7" = "The candidate has experience with team communication and collaboration, but lacks experience with Big Data Frameworks such as Hadoop, Spark, and Flink."
    "C8" = "The candidate has a strong technical background, relevant education, and some experience in software development, making them a good match for the role."
    "E8" = "Raleigh, North Carolina 27606 → United States (0 miles)"
    "G8" = "The candidate has some experience in software development, but lacks direct experience in search engine development and has limited overall work experience."
    "I8" = "The candidate has a strong educational background in computer science, making them a good match for the role."
    "K8" = "The candidate has a strong technical skillset, but lacks experience with some of the required programming languages and has limited experience with machine learning."
    "M8" = "The candidate speaks English, which is the required language for the role."
    "O8" = "The candidate has some experience with team communication and collaboration, but lacks direct experience with Big Data Frameworks."
    "C9" = "The candidate has a strong background in computer science and software engineering, with relevant experience in programming languages and algorithm development. However, there are some gaps in experience and skills required for the job, such as machine learning and deep learning."
    "E9" = "Richmond, North Carolina 27606 → United States (0 miles)"
    "G9" = "The candidate has relevant experience in software development, but lacks experience in some of the required programming languages."
    "I9" = "The candidate has a bachelor's degree in computer science, which meets the job requirement."
    "K9" = "The candidate has some of the required skills, but lacks experience in some of the programming languages."
    "M9" = "The candidate speaks English, which is the required language for the role."
    "O9" = "The candidate has experience with team communication and collaboration, but lacks experience with Big Data Frameworks."
    "C10" = "The candidate's background is in construction management, which does not align with the job requirements for a Backend Software Engineer."
    "E10" = "Callahan, FL 32011 → Culver City, CA (2500 miles)"
    "G10" = "The candidate has no experience in software development, which is a critical requirement for the job."
    "I10" = "The candidate's education in Civil Engineering does not align with the job requirements for a Backend Software Engineer"
    "K10" = "The candidate does not have the required programming skills, such as proficiency in languages like C/C++, Python, Java, or Golang."
    "M10" = "The candidate does not have any of the required programming languages."
    "O10" = "The candidate does not have any experience in the required areas, such as search, recommendation, ranking, ads, etc."
    "C11" = "The candidate's background is in environmental consulting and does not align with the job requirements for a backend software engineer."
    "E11" = "New York City → Culver City, CA (2500 miles)"
    "G11" = "The candidate has no experience in software development and does not meet the experience requirements for the job."
    "I11" = "The candidate has a degree in Geology, but it does not meet the education requirements for the job."
    "K11" = "The candidate does not have the required technical skills for the job."
    "M11" = "The candidate does not have any of the required programming languages."
    "O11" = "The candidate does not meet the other requirements for the job, including experience with machine learning and large-scale systems."
}

foreach ($addr in $notes.Keys) {
    $ws.Range($addr).Value = $notes[$addr]
}

# Header + data rows all grew to a taller 30pt row height.
$ws.Rows("1:11").RowHeight = 30

# Final selection left on F12 (scrolled so column E is no longer pinned as the
# left-most visible column).
$ws.Range("F12").Select() | Out-Null

Write-Output "done"
